$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2023-2024")

$ws.Range("A25").Value = Get-Date -Year 2024 -Month 3 -Day 11 -Hour 0 -Minute 0 -Second 0
$ws.Range("B25").Value = "FSIL"
$ws.Range("C25").Value = "TP"
$ws.Range("E25").Value = "x"
$ws.Range("G25").Value = "Suite questionScore : fix #2, fix #3. RAF : correction application."

$ws.Rows.Item(25).RowHeight = 39.75
